# "added update player area after discard"
# Fill in the "player" (B column) claimed-points values that mirror the
# already-present grader values in column C, plus a few "x" / "tbd"
# markers, and add a new "Claimed Total" summary row (134).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- B column: plain numeric values (style already s="5" or s="11") ---
$numericCells = @{
  "B23" = 5
  "B32" = 1
  "B33" = 1
  "B34" = 1
  "B35" = 1
  "B36" = 1
  "B48" = 2
  "B50" = 1
  "B57" = 1
  "B62" = 1
  "B63" = 1
  "B64" = 1
  "B93" = 0.5
  "B94" = 0.5
  "B95" = 0.5
  "B96" = 0.5
  "B97" = 0.5
  "B98" = 0.5
  "B99" = 1
  "B102" = 1
  "B103" = 1
  "B105" = 1
  "B106" = 0.5
  "B110" = 1
  "B111" = 0.5
  "B113" = 0
  "B114" = 1
  "B115" = 1
  "B118" = 1
  "B119" = 1
  "B124" = 0
}

# --- B column: cells marked with "x" instead of a number ---
$xCells = @("B51","B100","B101","B104","B107","B109","B112","B116","B117")

# --- D column: cells marked "tbd" (added alongside some of the rows above) ---
$tbdCells = @("D97","D98","D102","D103","D105","D106","D110","D111","D114","D115","D118","D119")

# Order matters for shared-string allocation: "x" first, then "tbd", then
# "Claimed Total" -- matches how the strings were newly introduced.
foreach ($addr in $xCells) {
  $ws.Range($addr).Value = "x"
}

foreach ($addr in $tbdCells) {
  $ws.Range($addr).Value = "tbd"
}

foreach ($addr in $numericCells.Keys) {
  $ws.Range($addr).Value = $numericCells[$addr]
}

# --- New summary row 134: Claimed Total ---
$ws.Range("A134").Value = "Claimed Total"
$ws.Range("B134").Formula = "=SUM(B32:B131)"
$ws.Range("C134").Formula = "=B134/(100-8)"

# --- Update the view so the selection matches the post-edit state ---
$ws.Range("D136").Select()
